$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("C17").Value = 45
$ws.Range("G17").Value = 21

# Row 18 (text-formatted numeric cell)
$ws.Range("C18").Value = "'25"

# Row 19 (text-formatted numeric cell)
$ws.Range("D19").Value = "'135"

# Row 20
$ws.Range("D20").Value = "22.2% der Karten"

# Row 21
$ws.Range("A21").Value = 39
$ws.Range("B21").Value = 96

# Row 27
$ws.Range("C27").Value = 32

# Row 28
$ws.Range("B28").Value = "Finanzfeedback für FS ⭐️"
$ws.Range("C28").Value = 16

# Row 29
$ws.Range("B29").Value = "Organigram aktualisieren ⭐️"
$ws.Range("C29").Value = 14

# Row 30
$ws.Range("B30").Value = "Ressorttreffen Projekte 29.06.2020 ⭐️"
$ws.Range("C30").Value = 14

# Row 31
$ws.Range("B31").Value = "Workshop Nutrición Escolar - Ernährungssicherheit ⭐️"
$ws.Range("C31").Value = 14

# Row 34
$ws.Range("F34").Value = 4
$ws.Range("G34").Value = "(8.9%)"

# Row 35
$ws.Range("C35").Value = 23

# Row 36
$ws.Range("C36").Value = 22

# Row 37
$ws.Range("C37").Value = 10

# Row 39
$ws.Range("B39").Value = "Laura Coordt"
$ws.Range("C39").Value = 2
$ws.Range("F39").Value = 19

# Row 40
$ws.Range("F40").Value = 18

# Row 43
$ws.Range("E43").Value = "Laura Coordt"
$ws.Range("F43").Value = 1

# Row 50
$ws.Range("B50").Value = "Valentin Buchenroth"
$ws.Range("F50").Value = 4

# Row 51
$ws.Range("B51").Value = "Jonas Ullmann"
$ws.Range("C51").Value = 1
$ws.Range("F51").Value = 4

# Row 52 - clear B52 (name removed), update counts
$ws.Range("B52").ClearContents()
$ws.Range("C52").Value = 0
$ws.Range("F52").Value = 3

# Row 53
$ws.Range("E53").Value = "Theresa Rinnert"

# Row 54
$ws.Range("E54").Value = "Lucia Irene Trepp"

# Row 60
$ws.Range("B60").Value = 7

# Row 61
$ws.Range("G61").Value = 14

# Row 62
$ws.Range("A62").Value = "Marie-Sophie Braun"

# Row 63 - add A63, update B63
$ws.Range("A63").Value = "Jonas Ullmann"
$ws.Range("B63").Value = 1

# Row 70
$ws.Range("B70").Value = 45
$ws.Range("G70").Value = 14

# Row 71
$ws.Range("B71").Value = 30
$ws.Range("G71").Value = 12

# Row 74
$ws.Range("G74").Value = 1
